# researcher_metadata_figshare.xlsx
# Add stage flag to get_<xxx>_info.py: authors and categories have
# different ids for production vs stage environments, so refresh the
# scraped author/category id values against the stage ids and drop the
# now-redundant second "authors id" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second "authors" id row (authors / 2 / id / 9314359) duplicated
# information no longer needed now that ids are environment-specific -
# remove it; rows below shift up by one.
$ws.Rows.Item(5).Delete()

# authors / 1 / id -> stage id. Pick up the "explicit black font" style
# already used for numeric id cells elsewhere in the sheet (e.g. D9:D12)
# instead of the default Attr_value style.
$ws.Range("D3").Value = 2933718
$ws.Range("D3").Font.Color = $ws.Range("D10").Font.Color()

# categories / 1 / categories -> stage id, and match styling used by the
# other category-id cells.
$ws.Range("D9").Value = 25718
$ws.Range("D9").Font.Color = $ws.Range("D10").Font.Color()

# categories / 1 / categories_by_source_id -> stage id
$ws.Range("D10").Value = 310112

# categories / 2 / categories -> stage id
$ws.Range("D11").Value = 26104

# categories / 2 / categories_by_source_id -> stage id
$ws.Range("D12").Value = 320999

# Nudge column D back to the sheet's default width (cosmetic leftover
# from editing the id values).
$ws.Columns.Item(4).ColumnWidth = 10

# Restore selection to where the author left off.
$ws.Range("I23").Select()
